$wb = $excel.ActiveWorkbook

# ALC!row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 350
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 350
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = -180
$ws.Range("N12").Value = -540

# ALC!row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 751.1132
$ws.Range("J17").Value = 764.8823
$ws.Range("L17").Value = 2294.6469
$ws.Range("N17").Value = -2630.6469

# ALC!row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 121.14286
$ws.Range("I33").Value = 89.916664
$ws.Range("J33").Value = 308.5
$ws.Range("K33").Value = 89.916664
$ws.Range("L33").Value = 308.5
$ws.Range("M33").Value = 139.083336
$ws.Range("N33").Value = -766.5

# ALC!row 123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 97208
$ws.Range("J123").Value = 97208
$ws.Range("L123").Value = 97208
$ws.Range("N123").Value = -107008

# ALC!row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5594342
$ws.Range("I138").Value = 948760.0600000001
$ws.Range("J138").Value = 8066990.5
$ws.Range("K138").Value = 2846280.18
$ws.Range("L138").Value = 24200971.5
$ws.Range("M138").Value = -2841140.18
$ws.Range("N138").Value = -24211251.5

# ARM!row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1940.381
$ws.Range("I45").Value = 1334.25
$ws.Range("J45").Value = 3880
$ws.Range("K45").Value = 1334.25
$ws.Range("L45").Value = 3880
$ws.Range("M45").Value = -957.25
$ws.Range("N45").Value = -4634

# ARM!row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1361.2858
$ws.Range("I102").Value = 1076
$ws.Range("K102").Value = 1076
$ws.Range("M102").Value = 546

# ARM!row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1759.6364
$ws.Range("I122").Value = 1178
$ws.Range("K122").Value = 3534
$ws.Range("M122").Value = -1084

# ARM!row 123
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 33607.332
$ws.Range("J123").Value = 33607.332
$ws.Range("L123").Value = 33607.332
$ws.Range("N123").Value = -43407.332

# ARM!row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1635.5143
$ws.Range("I132").Value = 1302.4333
$ws.Range("J132").Value = 3634
$ws.Range("K132").Value = 3907.2999
$ws.Range("L132").Value = 10902
$ws.Range("M132").Value = -1377.2999
$ws.Range("N132").Value = -15962

# ARM!row 135
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 34500
$ws.Range("J135").Value = 34500
$ws.Range("L135").Value = 34500
$ws.Range("N135").Value = -44640

# ARM!row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 64345.43
$ws.Range("J139").Value = 64345.43
$ws.Range("L139").Value = 64345.43
$ws.Range("N139").Value = -74625.42999999999

# BSM!row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11474.546
$ws.Range("I86").Value = 1952.1666
$ws.Range("J86").Value = 22901.4
$ws.Range("K86").Value = 1952.1666
$ws.Range("L86").Value = 22901.4
$ws.Range("M86").Value = -829.1666
$ws.Range("N86").Value = -25147.4

# BSM!row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 11474.546
$ws.Range("I89").Value = 1952.1666
$ws.Range("J89").Value = 22901.4
$ws.Range("K89").Value = 9760.833000000001
$ws.Range("L89").Value = 114507
$ws.Range("M89").Value = -4144.833000000001
$ws.Range("N89").Value = -125739

# BSM!row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 998.0345
$ws.Range("I107").Value = 880.381
$ws.Range("J107").Value = 1306.875
$ws.Range("K107").Value = 880.381
$ws.Range("L107").Value = 1306.875
$ws.Range("M107").Value = 1039.619
$ws.Range("N107").Value = -5146.875

# BSM!row 138
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 42166.668
$ws.Range("J138").Value = 42166.668
$ws.Range("L138").Value = 42166.668
$ws.Range("N138").Value = -52446.668

# CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1717.1187
$ws.Range("I31").Value = 1216.3334
$ws.Range("J31").Value = 2139.6562
$ws.Range("K31").Value = 1216.3334
$ws.Range("L31").Value = 2139.6562
$ws.Range("M31").Value = -921.3334
$ws.Range("N31").Value = -2729.6562

# CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1717.1187
$ws.Range("I34").Value = 1216.3334
$ws.Range("J34").Value = 2139.6562
$ws.Range("K34").Value = 1216.3334
$ws.Range("L34").Value = 2139.6562
$ws.Range("M34").Value = -1014.3334
$ws.Range("N34").Value = -2543.6562

# CRP!row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 718.3333
$ws.Range("I107").Value = 710
$ws.Range("J107").Value = 729.44446
$ws.Range("K107").Value = 710
$ws.Range("L107").Value = 729.44446
$ws.Range("M107").Value = 1210
$ws.Range("N107").Value = -4569.44446

# CUL!row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 909.36365
$ws.Range("I80").Value = 651
$ws.Range("K80").Value = 1953
$ws.Range("M80").Value = -1017

# CUL!row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 909.36365
$ws.Range("I83").Value = 651
$ws.Range("K83").Value = 5859
$ws.Range("M83").Value = -1179

# CUL!row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1218.8182
$ws.Range("I129").Value = 793.3333
$ws.Range("J129").Value = 1378.375
$ws.Range("K129").Value = 2379.9999
$ws.Range("L129").Value = 4135.125
$ws.Range("M129").Value = 2620.0001
$ws.Range("N129").Value = -14135.125

# CUL!row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 11355
$ws.Range("I136").Value = 2015
$ws.Range("J136").Value = 13430.556
$ws.Range("K136").Value = 6045
$ws.Range("L136").Value = 40291.66800000001
$ws.Range("M136").Value = -945
$ws.Range("N136").Value = -50491.66800000001

# GSM!row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1852987.4
$ws.Range("I122").Value = 2778481
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 8335443
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -8332993
$ws.Range("N122").Value = -10900

# GSM!row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1987.9354
$ws.Range("I126").Value = 1344
$ws.Range("K126").Value = 4032
$ws.Range("M126").Value = -1562

# GSM!row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2360.1355
$ws.Range("I132").Value = 1815.2941
$ws.Range("J132").Value = 3101.12
$ws.Range("K132").Value = 5445.8823
$ws.Range("L132").Value = 9303.360000000001
$ws.Range("M132").Value = -2915.8823
$ws.Range("N132").Value = -14363.36

# LTW!row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 800
$ws.Range("I16").Value = 800
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -630

# LTW!row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 666.6667
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 500
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 500
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -876

# LTW!row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 56756.223
$ws.Range("I82").Value = 100960.6
$ws.Range("J82").Value = 1500.75
$ws.Range("K82").Value = 100960.6
$ws.Range("L82").Value = 1500.75
$ws.Range("M82").Value = -100599.6
$ws.Range("N82").Value = -2222.75

# LTW!row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 56756.223
$ws.Range("I85").Value = 100960.6
$ws.Range("J85").Value = 1500.75
$ws.Range("K85").Value = 100960.6
$ws.Range("L85").Value = 1500.75
$ws.Range("M85").Value = -99712.60000000001
$ws.Range("N85").Value = -3996.75

# LTW!row 123
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 39000
$ws.Range("J123").Value = 39000
$ws.Range("L123").Value = 39000
$ws.Range("N123").Value = -48800

# WVR!row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1124.6666
$ws.Range("I113").Value = 1315.3
$ws.Range("K113").Value = 3945.9
$ws.Range("M113").Value = -1775.9

# WVR!row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 33710
$ws.Range("J123").Value = 33710
$ws.Range("L123").Value = 33710
$ws.Range("N123").Value = -43510

# WVR!row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 16668389
$ws.Range("I136").Value = 25641880
$ws.Range("K136").Value = 76925640
$ws.Range("M136").Value = -76923090
